# Append a new data row (row 3) to the "Artfynd" sheet, mirroring the
# structure of the existing row 2. Numeric-looking text fields are written
# with a leading apostrophe so Excel keeps them as text instead of
# auto-converting them to numbers/dates, matching the source data's typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Taxon / record identifiers ---------------------------------------
$ws.Range("A3").Value = 3060925
$ws.Range("B3").Value = 96334
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 220787
$ws.Range("F3").Value = "Knärot"
$ws.Range("G3").Value = "Goodyera repens"
$ws.Range("H3").Value = "(L.) R. Br."
$ws.Range("I3").Value = "'280"
$ws.Range("J3").Value = "plantor/tuvor"

# --- Location -----------------------------------------------------------
$ws.Range("P3").Value = "Kullen, 200 m ÖSÖ om, Vrm"
$ws.Range("Q3").Value = 363020.7647424658
$ws.Range("R3").Value = 6619760.180664059
$ws.Range("S3").Value = 50
$ws.Range("T3").Value = "Värmland"
$ws.Range("U3").Value = "Arvika"
$ws.Range("V3").Value = "Värmland"
$ws.Range("W3").Value = "Arvika"

# --- Dates / times (kept as text, like the source sheet) ---------------
$ws.Range("Y3").Value = "'2012-03-30"
$ws.Range("Z3").Value = "'00:00"
$ws.Range("AA3").Value = "'2012-03-30"
$ws.Range("AB3").Value = "'00:00"

# --- Comments / flags / habitat -----------------------------------------
$ws.Range("AC3").Value = "280 bladrosetter med 12 fröställningar, spridda från 668-632 till 770-669, längst i söder exponerade mot nyupptaget hygge, nyligen gallrat bestånd med mycket ris och körspår"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AH3").Value = "Granskog"
$ws.Range("AI3").Value = "nyligen gallrad blandbarrskog med gran och tall"

# AT3 is present but blank in the source row - force a text cell to exist
# with an empty value (mirrors the blank "Bestämningsår" column on row 2).
$ws.Range("AT3").NumberFormat = "@"
$ws.Range("AT3").Value = ""

# --- Reporter / observers ------------------------------------------------
$ws.Range("AW3").Value = "Per Larsson"
$ws.Range("AX3").Value = "Per Larsson"

# AY3 is likewise present but blank (mirrors the blank "Projektnamn" column).
$ws.Range("AY3").NumberFormat = "@"
$ws.Range("AY3").Value = ""
